$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 811d2a96-...md (row 3) and d46c28dd-...md (row 5)
$wsOverview.Range("G3").Value = "2016-09-06 10:21:34"
$wsOverview.Range("G5").Value = "2016-09-06 10:21:34"

# zh-cn sheet: row 3 (811d2a96) Priority "ht" -> "mt", Handoff/Handback datetimes
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-09-06 10:21:29"
$wsZhCn.Range("K3").Value = "2016-09-06 10:21:47"

# zh-cn sheet: row 5 (d46c28dd) same shared strings
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H5").Value = "2016-09-06 10:21:29"
$wsZhCn.Range("K5").Value = "2016-09-06 10:21:47"

# de-de sheet: row 3 (811d2a96) Priority "ht" -> "mt", Handoff datetime (shared w/ overview), Handback datetime
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-09-06 10:21:34"
$wsDeDe.Range("K3").Value = "2016-09-06 10:21:55"

# de-de sheet: row 5 (d46c28dd) same shared strings
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H5").Value = "2016-09-06 10:21:34"
$wsDeDe.Range("K5").Value = "2016-09-06 10:21:55"
